$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at row 377, shifting existing rows 377-441 down to 382-446
$ws.Range("377:381").Insert()

# Row 377
$ws.Range("A377").Value = 6
$ws.Range("B377").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C377").Value = "Metropolitana"
$ws.Range("D377").Value = 44476
$ws.Range("E377").Value = 13
$ws.Range("F377").Value = 100112021
$ws.Range("G377").Value = "Ají"
$ws.Range("H377").Value = "Americana (o)"
$ws.Range("I377").Value = "Primera"
$ws.Range("J377").Value = 80
$ws.Range("K377").Value = 90000
$ws.Range("L377").Value = 95000
$ws.Range("M377").Value = 92812
$ws.Range("N377").Value = "$/caja 25 kilos"
$ws.Range("O377").Value = "Provincia de Limarí"
$ws.Range("P377").Value = 3712
$ws.Range("Q377").Value = 25
$ws.Range("R377").Value = "Hortaliza"

# Row 378
$ws.Range("A378").Value = 6
$ws.Range("B378").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C378").Value = "Metropolitana"
$ws.Range("D378").Value = 44476
$ws.Range("E378").Value = 13
$ws.Range("F378").Value = 100112021
$ws.Range("G378").Value = "Ají"
$ws.Range("H378").Value = "Americana (o)"
$ws.Range("I378").Value = "Segunda"
$ws.Range("J378").Value = 26
$ws.Range("K378").Value = 75000
$ws.Range("L378").Value = 75000
$ws.Range("M378").Value = 75000
$ws.Range("N378").Value = "$/caja 25 kilos"
$ws.Range("O378").Value = "Provincia de Limarí"
$ws.Range("P378").Value = 3000
$ws.Range("Q378").Value = 25
$ws.Range("R378").Value = "Hortaliza"

# Row 379
$ws.Range("A379").Value = 6
$ws.Range("B379").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C379").Value = "Metropolitana"
$ws.Range("D379").Value = 44476
$ws.Range("E379").Value = 13
$ws.Range("F379").Value = 100112021
$ws.Range("G379").Value = "Ají"
$ws.Range("H379").Value = "Chilena(o)"
$ws.Range("I379").Value = "Primera"
$ws.Range("J379").Value = 14
$ws.Range("K379").Value = 14000
$ws.Range("L379").Value = 15000
$ws.Range("M379").Value = 14571
$ws.Range("N379").Value = "$/saco 25 kilos"
$ws.Range("O379").Value = "Provincia de Huasco"
$ws.Range("P379").Value = 583
$ws.Range("Q379").Value = 25
$ws.Range("R379").Value = "Hortaliza"

# Row 380
$ws.Range("A380").Value = 6
$ws.Range("B380").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C380").Value = "Metropolitana"
$ws.Range("D380").Value = 44476
$ws.Range("E380").Value = 13
$ws.Range("F380").Value = 100112021
$ws.Range("G380").Value = "Ají"
$ws.Range("H380").Value = "Chilena(o)"
$ws.Range("I380").Value = "Segunda"
$ws.Range("J380").Value = 5
$ws.Range("K380").Value = 110000
$ws.Range("L380").Value = 110000
$ws.Range("M380").Value = 110000
$ws.Range("N380").Value = "$/caja 25 kilos"
$ws.Range("O380").Value = "Provincia de Huasco"
$ws.Range("P380").Value = 4400
$ws.Range("Q380").Value = 25
$ws.Range("R380").Value = "Hortaliza"

# Row 381
$ws.Range("A381").Value = 6
$ws.Range("B381").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C381").Value = "Metropolitana"
$ws.Range("D381").Value = 44476
$ws.Range("E381").Value = 13
$ws.Range("F381").Value = 100112021
$ws.Range("G381").Value = "Ají"
$ws.Range("H381").Value = "Inferno"
$ws.Range("I381").Value = "Primera"
$ws.Range("J381").Value = 20
$ws.Range("K381").Value = 40000
$ws.Range("L381").Value = 45000
$ws.Range("M381").Value = 43000
$ws.Range("N381").Value = "$/caja 15 kilos"
$ws.Range("O381").Value = "Provincia de Huasco"
$ws.Range("P381").Value = 2867
$ws.Range("Q381").Value = 15
$ws.Range("R381").Value = "Hortaliza"
